# Atualização automática: 2025-08-25 21:00:26
#
# Refreshes the detection rows for the second sighting of PLACA_20250717165933
# on 2025-08-07/08: the First_Detection_Image file names, the First_Coords
# bounding boxes, and the First_Confidence scores for rows 16-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I16:J18 hold numeric-looking text ("643,531,686,575", "0.76", ...). Excel's
# default Value setter would auto-parse those as numbers (stripping the
# thousand-separator commas / trailing zeros), so mark the range as Text
# first to preserve them as literal strings, exactly like the source data.
$coordsAndConfidence = $ws.Range("I16:J18")
$coordsAndConfidence.NumberFormat = "@"

# Row 16
$ws.Range("D16").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I16").Value = "643,531,686,575"
$ws.Range("J16").Value = "0.76"

# Row 17
$ws.Range("D17").Value = "image_20250807111314_ppp0.jpg"
$ws.Range("I17").Value = "794,481,830,526"
$ws.Range("J17").Value = "0.72"

# Row 18
$ws.Range("D18").Value = "image_20250808100711_ppp0.jpg"
$ws.Range("I18").Value = "1182,409,1232,451"
$ws.Range("J18").Value = "0.75"

# Restore the default cell style now that the text values are locked in, so
# no stray number-format styling lingers on these cells.
$coordsAndConfidence.Style = "Normal"
